$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "123.123.123" rut sample value to "rut1"
$ws.Range("A9").Value = "rut1"

# Add a new row 10: a test rut "rut2", same date as other rows, and a score
$ws.Range("A10").Value = "rut2"

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 45198

$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 1000

# Add an empty, underlined cell at C11 (test formatting before the exam)
$ws.Range("C11").Font.Underline = $true

$ws.Range("A11").Select()
